$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Control 39
$ws.Range("D2").Value = [double]"0.9999999999287728"
$ws.Range("E2").Value = [double]"0.9999999999287728"

# Row 3 - Control 17
$ws.Range("D3").Value = [double]"6.397240495861397E-16"
$ws.Range("E3").Value = [double]"6.397240495861397E-16"

# Row 4 - Control 23
$ws.Range("D4").Value = [double]"0.00237195744260533"
$ws.Range("E4").Value = [double]"0.00237195744260533"

# Row 5 - Control 27
$ws.Range("D5").Value = [double]"1.950219491160004E-10"
$ws.Range("E5").Value = [double]"1.950219491160004E-10"

# Row 6 - Control 8
$ws.Range("D6").Value = [double]"8.887446939718487E-11"
$ws.Range("E6").Value = [double]"8.887446939718487E-11"

# Row 7 - MDD 27
$ws.Range("D7").Value = [double]"0.9999999999999836"
$ws.Range("E7").Value = [double]"1.643130076445232E-14"

# Row 8 - MDD 47
$ws.Range("D8").Value = [double]"8.411375117423039E-07"
$ws.Range("E8").Value = [double]"0.9999991588624882"

# Row 9 - MDD 13
$ws.Range("D9").Value = [double]"0.543898696947709"
$ws.Range("E9").Value = [double]"0.456101303052291"

# Row 11 - MDD 5
$ws.Range("D11").Value = [double]"0.9999835054572835"
$ws.Range("E11").Value = [double]"1.649454271646089E-05"
$ws.Range("F11").Value = [double]"3.796504259109497"
